# Add a new "2021" column (R) to the 9.2.1 table, mirroring the existing
# "2020" column (Q) for layout/formatting, then fill in the new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column Q's per-cell formatting into column R for the rows that carry
# data/formatting (header spacer row 2, year-header row 3, and the two data
# rows 4-5). Using Copy(destination) (rather than Copy()+PasteSpecial)
# preserves the source cell's formatting on the new cell.
$ws.Range("Q2").Copy($ws.Range("R2"))
$ws.Range("Q3").Copy($ws.Range("R3"))
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("Q5").Copy($ws.Range("R5"))

# Now set the new column's values: year 2021, and its two data points.
$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 13.5
$ws.Range("R5").Value = 15.1

# Match the author's final selection/active cell.
[void]$ws.Range("T3").Select()
